# additions to refinement scheme
# Add a new "Properties few" worksheet right after the existing "Properties"
# sheet, populate it with a small excerpt of the data (reusing the shared
# "failed" header string), and leave the selections where the author left
# them: P12 on "Properties" and A3 on the new "Properties few" sheet.

$wb = $excel.ActiveWorkbook
$props = $wb.Worksheets.Item("Properties")

$new = $wb.Worksheets.Add($null, $props)
$new.Name = "Properties few"

$new.Range("A1").Value = "failed"
$new.Range("A2").Value = 6000
$new.Range("A3").Value = 16000

$new.Range("A3").Select() | Out-Null
$props.Range("P12").Select() | Out-Null
